$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that moved from 2023-09-20
# (45189) to 2023-09-21 (45190) for every data row (2..270).
$ws.Range("C2:C270").Value = 45190
